$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as literal text in the source data (prices with
# dotted thousand separators etc.), so force Text number format before writing
# them to avoid Excel auto-converting them to numeric/date values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.987.10'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.848.63'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.012'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '309.74'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4778'
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3682'
$ws.Range("E8").Value = '  +2.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07233'
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9292'
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.73'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07743'
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.794.63'
$ws.Range("E13").Value = '  -2.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.345'
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.441'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.78'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008643'
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.014.94'
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.66'
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.921'
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.91'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  +1.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.002'
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.965'
$ws.Range("E29").Value = '  +1.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08884'
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.327'
$ws.Range("E31").Value = '  +5.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.172'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7421'
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.502'
$ws.Range("E34").Value = '  +1.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.751'
$ws.Range("E35").Value = '  -3.35%  '
$ws.Range("E36").Value = '  +3.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01963'
$ws.Range("E37").Value = '  +1.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05274'
$ws.Range("E38").Value = '  +2.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.978'
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5220'
$ws.Range("E40").Value = '  +3.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.990'
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1513'
$ws.Range("E42").Value = '  +0.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.231'
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.61'
$ws.Range("E44").Value = '  +5.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4752'
$ws.Range("E45").Value = '  +2.08%  '
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.73'
$ws.Range("E47").Value = '  +3.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.610'
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '65.56'
$ws.Range("E49").Value = '  +2.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06066'
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8884'
$ws.Range("E51").Value = '  +4.29%  '
